$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Version: 1.0.1 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text fix (applies to all 4 occurrences: B8, B17, B25, B33)
$ws.Range("B8").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B17").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B25").Value = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B33").Value = "O usuário devidamente autenticado e na tela inicial do sistema."

# TC1 expected results text fix (accents + spacing)
$ws.Range("D11").Value = "SYSTEM Exibe a lista de diárias (solicitações) aptas para pagamento ordenado pelo número da diária em ordem crescente. Exibe esta lista de diárias também ordenada pela data de chegada da solicitação na fase de liquidação (após registrar o empenho)."

# TC2 expected results: add trailing period
$ws.Range("D20").Value = "SYSTEM Apresenta a tela de Registrar Liquidações."

# TC3 expected results: remove duplicated "o nome"
$ws.Range("D28").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde deverá constar o nome do usuário logado (que se atribuiu como responsável pela liquidação) no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

# TC4 expected results: add trailing period
$ws.Range("D36").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."
